# Edit script implementing the commit:
# "Vendor Addition and Client Orders are collect, Need to run expense now."
#
# Summary of changes:
#  1. vendor_inventory: bump quantities for Kates Car/Tires & Kates Car/Gas,
#     add two new vendor-item rows (Test/Test, Hell/Souls).
#  2. clients: stray header labels typed into F1/G1 (Column1/Column2).
#  3. Two new worksheets ("Test", then "Hell") added right after
#     vendor_inventory - per-vendor copies of vendor_inventory, trimmed down
#     to that vendor's own row.
#  4. vendors: two new vendor names (Hell, Test) inserted near the top of the
#     list.
#
# (Order below mirrors the shared-string allocation order recoverable from
# the saved file: Test/Hell/Souls first, then Column1/Column2.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. vendor_inventory
# ---------------------------------------------------------------------------
$vi = $wb.Worksheets.Item("vendor_inventory")

$vi.Range("E2").Value = 8
$vi.Range("E3").Value = 6

$vi.Range("A11").Value = "Test"
$vi.Range("B11").Value = "Test"
$vi.Range("C11").Value = 5.99
$vi.Range("D11").Value = 0
$vi.Range("E11").Value = 3

$vi.Range("A12").Value = "Hell"
$vi.Range("B12").Value = "Souls"
$vi.Range("C12").Value = 6.99
$vi.Range("D12").Value = 0
$vi.Range("E12").Value = 3

# ---------------------------------------------------------------------------
# 2. clients
# ---------------------------------------------------------------------------
$cl = $wb.Worksheets.Item("clients")
$cl.Range("F1").Value = "Column1"
$cl.Range("G1").Value = "Column2"

# ---------------------------------------------------------------------------
# 3. New worksheets "Test" and "Hell" - copied from vendor_inventory as it
#    stood (with the two new rows already on it), then trimmed down to the
#    single matching vendor row. "Test" is copied first, "Hell" second so
#    that "Hell" lands immediately after vendor_inventory (ahead of "Test")
#    while still picking up the earlier internal sheetId.
# ---------------------------------------------------------------------------
$vi.Copy([System.Reflection.Missing]::Value, $vi)
$testSheet = $wb.Worksheets.Item("vendor_inventory (2)")
$testSheet.Name = "Test"

$testSheet.Range("A2").Value = "Test"
$testSheet.Range("B2").Value = "Test"
$testSheet.Range("C2").Value = 5.99
$testSheet.Range("D2").Value = 0
$testSheet.Range("E2").Value = 3
$testSheet.Range("A3:E10").ClearContents()
$testSheet.Range("A11:E12").EntireRow.Delete()

$vi.Copy([System.Reflection.Missing]::Value, $vi)
$hellSheet = $wb.Worksheets.Item("vendor_inventory (2)")
$hellSheet.Name = "Hell"

$hellSheet.Range("A2").Value = "Hell"
$hellSheet.Range("B2").Value = "Souls"
$hellSheet.Range("C2").Value = 6.99
$hellSheet.Range("D2").Value = 0
$hellSheet.Range("E2").Value = 3
$hellSheet.Range("A3:E12").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 4. vendors
# ---------------------------------------------------------------------------
$ve = $wb.Worksheets.Item("vendors")
$ve.Rows.Item(2).Insert()
$ve.Rows.Item(2).Insert()
$ve.Range("A2").Value = "Hell"
$ve.Range("A3").Value = "Test"
$ve.Range("A1:A4").Font.Bold = $false
$ve.Range("A16").EntireRow.Delete()
$ve.Rows.Item(19).RowHeight = 12.75
$ve.Rows.Item(20).RowHeight = 12.75
$ve.Rows.Item(21).RowHeight = 12.75

$cl.Range("E1").Select()

Write-Host "Edit complete"
